# Update "想去人数" (F column) values on sheet "展览" (rows 2-9)
# and on sheet "全部类型" (rows 2-6, 8-10), matching the new data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsAll = $wb.Worksheets.Item("全部类型")

# New F-column values keyed by row number, in the order they appear on "展览"
$exhibitValues = @{
    2 = 279
    3 = 175
    4 = 2078
    5 = 1658
    6 = 299
    7 = 85
    8 = 683
    9 = 151
}

foreach ($row in $exhibitValues.Keys) {
    $wsExhibit.Range("F$row").Value = $exhibitValues[$row]
}

# "全部类型" has the same events except row 7 (which is a "演出" entry, unaffected)
$allValues = @{
    2 = 279
    3 = 175
    4 = 2078
    5 = 1658
    6 = 299
    8 = 85
    9 = 683
    10 = 151
}

foreach ($row in $allValues.Keys) {
    $wsAll.Range("F$row").Value = $allValues[$row]
}
